$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.794.90'
$ws.Range('E2').Value = '  +3.63%  '
$ws.Range('D3').Value = '1.863.16'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '268.34'
$ws.Range('E5').Value = '  -3.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5293'
$ws.Range('E7').Value = '  +3.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3346'
$ws.Range('E8').Value = '  -3.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06789'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.65'
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7814'
$ws.Range('E11').Value = '  -3.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07754'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').Value = '1.880.44'
$ws.Range('E13').Value = '  +3.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.95'
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.108'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.35'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007988'
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '26.811.14'
$ws.Range('E20').Value = '  +3.58%  '
$ws.Range('D21').Value = '2.096.60'
$ws.Range('E21').Value = '  +2.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.651'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.854'
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.046'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.408'
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.50'
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.664'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.73'
$ws.Range('E29').Value = '  +3.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.281'
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.265'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08829'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04936'
$ws.Range('E33').Value = '  +1.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.154'
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.881'
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7198'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.186'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01830'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.292'
$ws.Range('E39').Value = '  -4.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5021'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '115.53'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9154'
$ws.Range('E42').Value = '  -3.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.115'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.925'
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4355'
$ws.Range('E46').Value = '  -2.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1314'
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.285'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05937'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.81'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '60.78'
$ws.Range('E51').Value = '  +0.68%  '
